$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A34:H34").Copy()
$ws.Range("A35:H35").Insert(-4121)
Write-Host "done"
